$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: new part "crossmember" (id 2) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "crossmember"
$ws.Range("C3").Value = 53489514
$ws.Range("D3").Value = 5980
$ws.Range("E3").Value = 45957.62032447917
$ws.Range("F3").Value = "dados/peca_2/txt"
$ws.Range("G3").Value = "Ativa"

# --- Row 4: new part "chifre" (id 3), PartNumber/Modelo recorded as text ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "chifre"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "53489571"
$ws.Range("C4").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5980"
$ws.Range("D4").ClearFormats()

$ws.Range("E2").Copy($ws.Range("E4"))
$ws.Range("E4").Value = 45957.62969117925

$ws.Range("F4").Value = "dados/peca_3/txt"
$ws.Range("G4").Value = "Ativa"
